# Weekly price-data update: a new record (week of 2022-08-09) is inserted
# at row 338 of the "Hortaliza, Feria Lagunitas de Puerto Montt - Zanahoria"
# sheet, pushing all subsequent records down by one row (old row 414 -> 415).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 338:414 down to 339:415, inserting a blank row at 338
# (mirrors Excel's Rows(...).Insert, including the dimension/ref update).
$ws.Rows("338:338").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A338").Value = 4
$ws.Range("B338").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C338").Value = "Los Lagos"
$ws.Range("D338").Value = 44782
$ws.Range("E338").Value = 10
$ws.Range("F338").Value = 100114013
$ws.Range("G338").Value = "Zanahoria"
$ws.Range("H338").Value = "Sin especificar"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 900
$ws.Range("K338").Value = 10000
$ws.Range("L338").Value = 10500
$ws.Range("M338").Value = 10250
$ws.Range("N338").Value = "`$/saco 20 kilos"
$ws.Range("O338").Value = "Provincia de Llanquihue"
$ws.Range("P338").Value = 512
$ws.Range("Q338").Value = 20
$ws.Range("R338").Value = "Hortaliza"
